$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray one-off formula in G33 (=30-13). Clearing it also causes
# the sheet's used-range dimension and each row's "spans" attribute to be
# recalculated automatically (G column no longer used => A1:E50).
$ws.Range("G33").ClearContents() | Out-Null

# Fill in the two new timesheet entries that were previously blank
# placeholder rows (45 & 46) - Date 12/3/2017, with start/end times.
$ws.Range("A45").Value = 43072
$ws.Range("B45").Value = 0.60416666666666663
$ws.Range("C45").Value = 0.73958333333333337

$ws.Range("A46").Value = 43072
$ws.Range("B46").Value = 0.77083333333333337
$ws.Range("C46").Value = 0.79166666666666663

# Restore the cursor/selection to what the author left it at (best-effort;
# the headless host does not model scroll position / topLeftCell).
$ws.Range("F40").Select() | Out-Null

Write-Output "edit applied"
